$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style used to strip the "quote prefix" (text-stored-as-number) marker
# that Excel applies automatically when a numeric-looking string is force-typed as text.
$normalStyle = $ws.Range("A1").Style

$ws.Range('D2').Value = "'26.802.32"
$ws.Range('D2').Style = $normalStyle
$ws.Range('E2').Value = '  -1.77%  '

$ws.Range('D3').Value = "'1.869.01"
$ws.Range('D3').Style = $normalStyle
$ws.Range('E3').Value = '  -2.07%  '

$ws.Range('D4').Value = "'1.000"
$ws.Range('D4').Style = $normalStyle

$ws.Range('D5').Value = "'300.00"
$ws.Range('D5').Style = $normalStyle
$ws.Range('E5').Value = '  -2.51%  '

$ws.Range('D6').Value = "'1.0000"
$ws.Range('D6').Style = $normalStyle
$ws.Range('E6').Value = '  -0.10%  '

$ws.Range('E7').Value = '  +1.69%  '

$ws.Range('D8').Value = "'0.3738"
$ws.Range('D8').Style = $normalStyle
$ws.Range('E8').Value = '  -2.18%  '

$ws.Range('D9').Value = "'0.07142"
$ws.Range('D9').Style = $normalStyle
$ws.Range('E9').Value = '  -2.23%  '

$ws.Range('E10').Value = '  -2.48%  '

$ws.Range('D11').Value = "'0.8870"
$ws.Range('D11').Style = $normalStyle
$ws.Range('E11').Value = '  -1.70%  '

$ws.Range('D12').Value = "'0.08138"
$ws.Range('D12').Style = $normalStyle
$ws.Range('E12').Value = '  -0.73%  '

$ws.Range('D13').Value = "'1.905.75"
$ws.Range('D13').Style = $normalStyle
$ws.Range('E13').Value = '  +44.64%  '

$ws.Range('D14').Value = "'92.46"
$ws.Range('D14').Style = $normalStyle
$ws.Range('E14').Value = '  -3.51%  '

$ws.Range('E15').Value = '  -1.24%  '

$ws.Range('D16').Value = "'1.001"
$ws.Range('D16').Style = $normalStyle
$ws.Range('E16').Value = '  +0.00%  '

$ws.Range('D17').Value = "'14.83"
$ws.Range('D17').Style = $normalStyle
$ws.Range('E17').Value = '  +0.27%  '

$ws.Range('D18').Value = "'0.000008481"
$ws.Range('D18').Style = $normalStyle
$ws.Range('E18').Value = '  -1.89%  '

$ws.Range('E19').Value = '  -0.10%  '

$ws.Range('D20').Value = "'26.825.61"
$ws.Range('D20').Style = $normalStyle
$ws.Range('E20').Value = '  -1.80%  '

$ws.Range('E21').Value = '  -2.03%  '

$ws.Range('D22').Value = "'10.64"
$ws.Range('D22').Style = $normalStyle
$ws.Range('E22').Value = '  -1.65%  '

$ws.Range('D23').Value = "'6.375"
$ws.Range('D23').Style = $normalStyle
$ws.Range('E23').Value = '  -2.23%  '

$ws.Range('D24').Value = "'2.283"
$ws.Range('D24').Style = $normalStyle
$ws.Range('E24').Value = '  -0.59%  '

$ws.Range('D25').Value = "'146.04"
$ws.Range('D25').Style = $normalStyle
$ws.Range('E25').Value = '  -2.56%  '

$ws.Range('D26').Value = "'1.742"
$ws.Range('D26').Style = $normalStyle
$ws.Range('E26').Value = '  +0.24%  '

$ws.Range('E27').Value = '  -1.55%  '

$ws.Range('D28').Value = "'113.72"
$ws.Range('D28').Style = $normalStyle
$ws.Range('E28').Value = '  -2.34%  '

$ws.Range('D29').Value = "'4.698"
$ws.Range('D29').Style = $normalStyle
$ws.Range('E29').Value = '  -2.68%  '

$ws.Range('D30').Value = "'4.623"
$ws.Range('D30').Style = $normalStyle

$ws.Range('D31').Value = "'0.09100"
$ws.Range('D31').Style = $normalStyle
$ws.Range('E31').Value = '  -1.85%  '

$ws.Range('D32').Value = "'0.8117"
$ws.Range('D32').Style = $normalStyle
$ws.Range('E32').Value = '  -2.98%  '

$ws.Range('D33').Value = "'0.05017"
$ws.Range('D33').Style = $normalStyle
$ws.Range('E33').Value = '  -1.12%  '

$ws.Range('D34').Value = "'1.170"
$ws.Range('D34').Style = $normalStyle
$ws.Range('E34').Value = '  -4.74%  '

$ws.Range('D35').Value = "'2.944"
$ws.Range('D35').Style = $normalStyle
$ws.Range('E35').Value = '  -2.01%  '

$ws.Range('D36').Value = "'0.6074"
$ws.Range('D36').Style = $normalStyle
$ws.Range('E36').Value = '  +5.55%  '

$ws.Range('D37').Value = "'2.661"
$ws.Range('D37').Style = $normalStyle
$ws.Range('E37').Value = '  -0.80%  '

$ws.Range('D38').Value = "'3.198"
$ws.Range('D38').Style = $normalStyle
$ws.Range('E38').Value = '  -4.57%  '

$ws.Range('D39').Value = "'0.01944"
$ws.Range('D39').Style = $normalStyle
$ws.Range('E39').Value = '  -3.06%  '

$ws.Range('E40').Value = '  -1.04%  '

$ws.Range('D41').Value = "'0.5300"
$ws.Range('D41').Style = $normalStyle
$ws.Range('E41').Value = '  +7.72%  '

$ws.Range('D42').Value = "'8.758"
$ws.Range('D42').Style = $normalStyle
$ws.Range('E42').Value = '  -6.43%  '

$ws.Range('D43').Value = "'6.472"
$ws.Range('D43').Style = $normalStyle
$ws.Range('E43').Value = '  -0.97%  '

$ws.Range('D44').Value = "'116.31"
$ws.Range('D44').Style = $normalStyle
$ws.Range('E44').Value = '  -0.44%  '

$ws.Range('D45').Value = "'0.1486"
$ws.Range('D45').Style = $normalStyle
$ws.Range('E45').Value = '  -2.43%  '

$ws.Range('E46').Value = '  -0.07%  '

$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'10.00"
$ws.Range('D47').Style = $normalStyle
$ws.Range('E47').Value = '  -1.63%  '

$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = "'1.644"
$ws.Range('D48').Style = $normalStyle
$ws.Range('E48').Value = '  +0.39%  '

$ws.Range('D49').Value = "'37.31"
$ws.Range('D49').Style = $normalStyle
$ws.Range('E49').Value = '  -3.97%  '

$ws.Range('D50').Value = "'0.06059"
$ws.Range('D50').Style = $normalStyle
$ws.Range('E50').Value = '  -2.14%  '

$ws.Range('D51').Value = "'62.11"
$ws.Range('D51').Style = $normalStyle
$ws.Range('E51').Value = '  -2.58%  '
